# Montenegro Prva Liga - base update (12-06-2024 23:38)
#
# The source feed re-ordered several fixtures that share the same kickoff
# date/time (the "Date" column, D, is identical within each pair). For each
# affected pair of adjacent data rows, every other column (match id, teams,
# score, odds, etc.) needs to be swapped between the two rows while each
# row keeps its own running index in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-number pairs (1-based worksheet rows) whose B:AD contents are swapped.
$rowPairs = @(
    @(14, 15),
    @(17, 18),
    @(48, 49),
    @(59, 60),
    @(100, 101),
    @(107, 108),
    @(163, 164)
)

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    $rangeA = $ws.Range("B$rowA`:AD$rowA")
    $rangeB = $ws.Range("B$rowB`:AD$rowB")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}
